$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.475.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.858.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4761"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3787"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07317"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9284"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07769"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.862.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.437"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.550"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008811"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.455.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.085"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.936"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.001"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.945"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08861"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("E32").Value = "  +2.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7494"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.576"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.702"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("E36").Value = "  +4.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.119"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5549"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05296"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.016"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.493"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4860"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.661"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9107"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.12%  "
